# Apply updates described by the diff for 2024-10-11 FlashScore workbook
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 7 (Tepatitlan de Morelos vs Tampico Madero): update odds ---
$ws.Range("H7").Value = 3.1
$ws.Range("J7").Value = 2.75
$ws.Range("K7").Value = 2.05
$ws.Range("L7").Value = 3.75
$ws.Range("T7").Value = 2.5
$ws.Range("U7").Value = 1.8
$ws.Range("V7").Value = 1.8
$ws.Range("W7").Value = 7.1
$ws.Range("X7").Value = 10.5
$ws.Range("AA7").Value = 19
$ws.Range("AB7").Value = 30
$ws.Range("AD7").Value = 6
$ws.Range("AE7").Value = 15
$ws.Range("AG7").Value = 8.25
$ws.Range("AH7").Value = 16
$ws.Range("AI7").Value = 11.5
$ws.Range("AK7").Value = 32
$ws.Range("AM7").Value = 700
$ws.Range("AO7").Value = 11.25
$ws.Range("AP7").Value = 19
$ws.Range("AR7").Value = 75
$ws.Range("AU7").Value = 6.9
$ws.Range("AX7").Value = 18
$ws.Range("AY7").Value = 25
$ws.Range("BB7").Value = 350

# --- Row 8: update odds ---
$ws.Range("G8").Value = 2.88
$ws.Range("I8").Value = 2.45
$ws.Range("J8").Value = 3.6
$ws.Range("L8").Value = 3.2
$ws.Range("Q8").Value = 2.2
$ws.Range("R8").Value = 1.65
$ws.Range("S8").Value = 1.5
$ws.Range("T8").Value = 2.5
$ws.Range("AA8").Value = 26
$ws.Range("AG8").Value = 7.5
$ws.Range("AH8").Value = 11
$ws.Range("AJ8").Value = 23
$ws.Range("AP8").Value = 29
$ws.Range("AT8").Value = 2.5
$ws.Range("AU8").Value = 8.5

# --- Insert a new row at position 12 (Penarol vs Maldonado, Uruguay) ---
# This shifts existing rows 12 (La Guaira vs Rayo Zuliano) and 13 (Zamora vs Monagas) down to 13 and 14
$ws.Rows.Item(12).Insert()

# --- Populate new row 12 ---
$ws.Range("A12").Value = "YNeegDWH"
$ws.Range("B12").NumberFormat = "@"
$ws.Range("B12").Value = "11/10/2024"
$ws.Range("C12").Value = "19:00"
$ws.Range("D12").Value = "URUGUAY - PRIMERA DIVISION"
$ws.Range("E12").Value = "Penarol"
$ws.Range("F12").Value = "Maldonado"
$ws.Range("G12").Value = 1.38
$ws.Range("H12").Value = 4.33
$ws.Range("I12").Value = 9.5
$ws.Range("J12").Value = 1.91
$ws.Range("K12").Value = 2.25
$ws.Range("L12").Value = 8.5
$ws.Range("M12").Value = 1.06
$ws.Range("N12").Value = 10
$ws.Range("O12").Value = 1.33
$ws.Range("P12").Value = 3.25
$ws.Range("Q12").Value = 2.08
$ws.Range("R12").Value = 1.73
$ws.Range("S12").Value = 1.44
$ws.Range("T12").Value = 2.63
$ws.Range("U12").Value = 2.5
$ws.Range("V12").Value = 1.5
$ws.Range("W12").Value = 5.5
$ws.Range("X12").Value = 5.5
$ws.Range("Y12").Value = 9
$ws.Range("Z12").Value = 8.5
$ws.Range("AA12").Value = 13
$ws.Range("AB12").Value = 41
$ws.Range("AC12").Value = 8.5
$ws.Range("AD12").Value = 8.5
$ws.Range("AE12").Value = 26
$ws.Range("AF12").Value = 101
$ws.Range("AG12").Value = 19
$ws.Range("AH12").Value = 41
$ws.Range("AI12").Value = 29
$ws.Range("AJ12").Value = 126
$ws.Range("AK12").Value = 81
$ws.Range("AL12").Value = 81
$ws.Range("AM12").Value = 201
$ws.Range("AN12").Value = 3.1
$ws.Range("AO12").Value = 6.5
$ws.Range("AP12").Value = 23
$ws.Range("AQ12").Value = 21
$ws.Range("AR12").Value = 51
$ws.Range("AS12").Value = 201
$ws.Range("AT12").Value = 2.63
$ws.Range("AU12").Value = 11
$ws.Range("AV12").Value = 81
$ws.Range("AW12").Value = 9
$ws.Range("AX12").Value = 41
$ws.Range("AY12").Value = 51
$ws.Range("AZ12").Value = 251
$ws.Range("BA12").Value = 301
$ws.Range("BB12").Value = 501
$ws.Range("BC12").Value = 51
$ws.Range("BD12").Value = 51

# --- Row 13 (was row 12, La Guaira vs Rayo Zuliano): update odds after shift ---
$ws.Range("I13").Value = 5
$ws.Range("J13").Value = 2.2
$ws.Range("K13").Value = 2.12
$ws.Range("L13").Value = 5.1
$ws.Range("S13").Value = 1.39
$ws.Range("T13").Value = 2.57
$ws.Range("W13").Value = 6.1
$ws.Range("AB13").Value = 30
$ws.Range("AC13").Value = 9
$ws.Range("AH13").Value = 28
$ws.Range("AI13").Value = 16
$ws.Range("AJ13").Value = 90
$ws.Range("AK13").Value = 55
$ws.Range("AL13").Value = 60
$ws.Range("AM13").Value = 800
$ws.Range("AO13").Value = 8
$ws.Range("AP13").Value = 17.5
$ws.Range("AQ13").Value = 26
$ws.Range("AR13").Value = 60
$ws.Range("AT13").Value = 2.55
$ws.Range("AU13").Value = 7.6
$ws.Range("AV13").Value = 75
$ws.Range("AW13").Value = 6.5
$ws.Range("AX13").Value = 28
$ws.Range("AZ13").Value = 175
$ws.Range("BA13").Value = 200
$ws.Range("BB13").Value = 450

# --- Row 14 (was row 13, Zamora vs Monagas): update odds after shift ---
$ws.Range("H14").Value = 2.92
$ws.Range("I14").Value = 2.35
$ws.Range("J14").Value = 3.6
$ws.Range("K14").Value = 1.98
$ws.Range("L14").Value = 3
$ws.Range("N14").Value = 7.1
$ws.Range("P14").Value = 2.62
$ws.Range("T14").Value = 2.45
$ws.Range("V14").Value = 1.83
$ws.Range("X14").Value = 16.5
$ws.Range("Z14").Value = 45
$ws.Range("AC14").Value = 7.7
$ws.Range("AD14").Value = 5.7
$ws.Range("AJ14").Value = 25
$ws.Range("AK14").Value = 22
$ws.Range("AO14").Value = 17
$ws.Range("AS14").Value = 250
$ws.Range("AU14").Value = 6.7
$ws.Range("AW14").Value = 4.2
$ws.Range("AX14").Value = 13

